$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69, pushing existing rows 69..80 down to 70..81
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new data record
$ws.Cells.Item(69, 1).Value = 5
$ws.Cells.Item(69, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(69, 3).Value = "Maule"
$ws.Cells.Item(69, 4).Value = 44543
$ws.Cells.Item(69, 5).Value = 7
$ws.Cells.Item(69, 6).Value = 100112022
$ws.Cells.Item(69, 7).Value = "Arveja Verde"
$ws.Cells.Item(69, 8).Value = "Sin especificar"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 200
$ws.Cells.Item(69, 11).Value = 13000
$ws.Cells.Item(69, 12).Value = 13000
$ws.Cells.Item(69, 13).Value = 13000
$ws.Cells.Item(69, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(69, 15).Value = "Carahue"
$ws.Cells.Item(69, 16).Value = 520
$ws.Cells.Item(69, 17).Value = 25
$ws.Cells.Item(69, 18).Value = "Hortaliza"
